# CIERRE 29 NOV 21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REMISIONES   NOVIEMBRE  2021 ")

# Orange fill color used to flag rows 27 & 28 (matches existing style fillId=5 / RGB FFC000)
$orange = 49407

# Row 27: mark as paid with "x", clear/fill D and E with orange highlight
$ws.Range("D27").Value = "x"
$ws.Range("D27").Interior.Color = $orange
$ws.Range("E27").Interior.Color = $orange

# Row 28: same treatment
$ws.Range("D28").Value = "x"
$ws.Range("D28").Interior.Color = $orange
$ws.Range("E28").Interior.Color = $orange

# Row 29
$ws.Range("A29").Value = 44516
$ws.Range("D29").Value = "COMERCIO CENTRAL "
$ws.Range("E29").Value = 4624

# Row 30
$ws.Range("A30").Value = 44516
$ws.Range("D30").Value = "OBRADOR"
$ws.Range("E30").Value = 3512

# Row 31
$ws.Range("A31").Value = 44517
$ws.Range("D31").Value = "COMERCIO CENTRAL "
$ws.Range("E31").Value = 178470

# Row 32
$ws.Range("A32").Value = 44517
$ws.Range("D32").Value = "COMERCIO CENTRAL "
$ws.Range("E32").Value = 62080

# Row 33
$ws.Range("A33").Value = 44517
$ws.Range("D33").Value = "COMERCIO CENTRAL "
$ws.Range("E33").Value = 101

# Row 34
$ws.Range("A34").Value = 44518
$ws.Range("D34").Value = "COMERCIO CENTRAL "
$ws.Range("E34").Value = 8588

# Row 35
$ws.Range("A35").Value = 44520
$ws.Range("D35").Value = "COMERCIO CENTRAL "
$ws.Range("E35").Value = 768

# Update the active selection to F35 as in the edited file
$ws.Activate()
$ws.Range("F35").Select()
